$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the affected rows, per repull of data.
$ws.Range("F2").Value = -9
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = -6
$ws.Range("F7").Value = -4
$ws.Range("F9").Value = -6
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = -2
$ws.Range("F13").Value = 0
